$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.570.12"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "3.115.07"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "616.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +4.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.392"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.828"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.99%  "
$ws.Range("D11").Value = "3.115.46"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").Value = "94.193.41"
$ws.Range("E14").Value = "  +2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "3.695.54"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "3.126.28"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "451.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("E25").Value = "  +5.35%  "
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "86.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  +3.72%  "
$ws.Range("D29").Value = "3.291.66"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.261"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.87%  "
$ws.Range("E32").Value = "  +8.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.127"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.455"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "479.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("E44").Value = "  -10.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.695"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("E51").Value = "  -3.00%  "
